$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "2025/12/03 17:00"
$ws.Range("B25").Value = "-"
$ws.Range("C25").Value = "-"
$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = "-"
$ws.Range("F25").Value = "-"
$ws.Range("G25").Value = "-"
